$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 556.9375
$ws.Range("I2").Value = 560.1429000000001
$ws.Range("K2").Value = 560.1429000000001
$ws.Range("M2").Value = -447.1429000000001
$ws.Range("H17").Value = 821.4167
$ws.Range("J17").Value = 1012.875
$ws.Range("L17").Value = 3038.625
$ws.Range("N17").Value = -3374.625
$ws.Range("H19").Value = 1545.8
$ws.Range("J19").Value = 1948.4
$ws.Range("L19").Value = 1948.4
$ws.Range("N19").Value = -2298.4
$ws.Range("H32").Value = 4891.95
$ws.Range("I32").Value = 5141.091
$ws.Range("J32").Value = 4587.4443
$ws.Range("K32").Value = 5141.091
$ws.Range("L32").Value = 4587.4443
$ws.Range("M32").Value = -4815.091
$ws.Range("N32").Value = -5239.4443
$ws.Range("H43").Value = 5974.3335
$ws.Range("I43").Value = 5999.5
$ws.Range("K43").Value = 5999.5
$ws.Range("M43").Value = -5930.5
$ws.Range("H62").Value = 24633.268
$ws.Range("I62").Value = 8321
$ws.Range("K62").Value = 8321
$ws.Range("M62").Value = -7697
$ws.Range("H65").Value = 24633.268
$ws.Range("I65").Value = 8321
$ws.Range("K65").Value = 41605
$ws.Range("M65").Value = -38485
$ws.Range("H76").Value = 3553.3
$ws.Range("I76").Value = 3191.625
$ws.Range("K76").Value = 3191.625
$ws.Range("M76").Value = -2876.625
$ws.Range("H79").Value = 3553.3
$ws.Range("I79").Value = 3191.625
$ws.Range("K79").Value = 3191.625
$ws.Range("M79").Value = -2099.625
$ws.Range("H98").Value = 41983.688
$ws.Range("I98").Value = 46562.082
$ws.Range("J98").Value = 28248.5
$ws.Range("K98").Value = 46562.082
$ws.Range("L98").Value = 28248.5
$ws.Range("M98").Value = -45064.082
$ws.Range("N98").Value = -31244.5
$ws.Range("H106").Value = 5038.375
$ws.Range("I106").Value = 5329.5713
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 5329.5713
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -4698.5713
$ws.Range("N106").Value = -4262
$ws.Range("H122").Value = 41983.688
$ws.Range("I122").Value = 46562.082
$ws.Range("J122").Value = 28248.5
$ws.Range("K122").Value = 139686.246
$ws.Range("L122").Value = 84745.5
$ws.Range("M122").Value = -137236.246
$ws.Range("N122").Value = -89645.5
$ws.Range("H132").Value = 3273.3635
$ws.Range("I132").Value = 3300.276
$ws.Range("K132").Value = 9900.828
$ws.Range("M132").Value = -7370.828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7157.864
$ws.Range("I2").Value = 8812.3125
$ws.Range("K2").Value = 8812.3125
$ws.Range("M2").Value = -8699.3125
$ws.Range("H116").Value = 7157.864
$ws.Range("I116").Value = 8812.3125
$ws.Range("K116").Value = 8812.3125
$ws.Range("M116").Value = -6518.3125
$ws.Range("H122").Value = 428712.9
$ws.Range("I122").Value = 3629.5
$ws.Range("K122").Value = 10888.5
$ws.Range("M122").Value = -8438.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7157.864
$ws.Range("I3").Value = 8812.3125
$ws.Range("K3").Value = 8812.3125
$ws.Range("M3").Value = -8698.3125
$ws.Range("H43").Value = 245196.12
$ws.Range("J43").Value = 245196.12
$ws.Range("L43").Value = 245196.12
$ws.Range("N43").Value = -245558.12

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1221.7142
$ws.Range("J122").Value = 1567.5714
$ws.Range("L122").Value = 4702.7142
$ws.Range("N122").Value = -9602.7142
$ws.Range("H134").Value = 3647.275
$ws.Range("I134").Value = 2993.889
$ws.Range("K134").Value = 8981.667000000001
$ws.Range("M134").Value = -6446.667000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 309.41934
$ws.Range("I17").Value = 94.625
$ws.Range("J17").Value = 1045.8572
$ws.Range("K17").Value = 283.875
$ws.Range("L17").Value = 3137.5716
$ws.Range("M17").Value = -114.875
$ws.Range("N17").Value = -3475.5716
$ws.Range("H34").Value = 1473097
$ws.Range("I34").Value = 2502275
$ws.Range("K34").Value = 7506825
$ws.Range("M34").Value = -7506741
$ws.Range("H51").Value = 1937.6666
$ws.Range("I51").Value = 922.2
$ws.Range("J51").Value = 3968.6
$ws.Range("K51").Value = 2766.6
$ws.Range("L51").Value = 11905.8
$ws.Range("M51").Value = -2306.6
$ws.Range("N51").Value = -12825.8
$ws.Range("H86").Value = 1102.8422
$ws.Range("I86").Value = 1359.8
$ws.Range("J86").Value = 1063.909
$ws.Range("K86").Value = 4079.4
$ws.Range("L86").Value = 3191.727
$ws.Range("M86").Value = -2893.4
$ws.Range("N86").Value = -5563.727000000001
$ws.Range("H89").Value = 1102.8422
$ws.Range("I89").Value = 1359.8
$ws.Range("J89").Value = 1063.909
$ws.Range("K89").Value = 12238.2
$ws.Range("L89").Value = 9575.181
$ws.Range("M89").Value = -6310.199999999999
$ws.Range("N89").Value = -21431.181
$ws.Range("H136").Value = 437.25
$ws.Range("I136").Value = 437.25
$ws.Range("K136").Value = 1311.75
$ws.Range("M136").Value = 3788.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 18827.824
$ws.Range("I122").Value = 15371.533
$ws.Range("K122").Value = 46114.599
$ws.Range("M122").Value = -43664.599
$ws.Range("H123").Value = 20785.572
$ws.Range("J123").Value = 20785.572
$ws.Range("L123").Value = 20785.572
$ws.Range("N123").Value = -25685.572
$ws.Range("H132").Value = 4513.4165
$ws.Range("I132").Value = 4469.1816
$ws.Range("K132").Value = 13407.5448
$ws.Range("M132").Value = -10877.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23205.773
$ws.Range("J7").Value = 7149
$ws.Range("L7").Value = 7149
$ws.Range("N7").Value = -7373
$ws.Range("H40").Value = 31005.375
$ws.Range("I40").Value = 43709.3
$ws.Range("J40").Value = 9832.166999999999
$ws.Range("K40").Value = 43709.3
$ws.Range("L40").Value = 9832.166999999999
$ws.Range("M40").Value = -43573.3
$ws.Range("N40").Value = -10104.167
$ws.Range("H46").Value = 4598.923
$ws.Range("J46").Value = 6535.875
$ws.Range("L46").Value = 6535.875
$ws.Range("N46").Value = -6911.875
$ws.Range("H55").Value = 579.2105
$ws.Range("I55").Value = 636.4
$ws.Range("J55").Value = 364.75
$ws.Range("K55").Value = 636.4
$ws.Range("L55").Value = 364.75
$ws.Range("M55").Value = -463.4
$ws.Range("N55").Value = -710.75
$ws.Range("H68").Value = 6889
$ws.Range("I68").Value = 6750.5
$ws.Range("J68").Value = 6999.8
$ws.Range("K68").Value = 6750.5
$ws.Range("L68").Value = 6999.8
$ws.Range("M68").Value = -6001.5
$ws.Range("N68").Value = -8497.799999999999
$ws.Range("H71").Value = 6889
$ws.Range("I71").Value = 6750.5
$ws.Range("J71").Value = 6999.8
$ws.Range("K71").Value = 33752.5
$ws.Range("L71").Value = 34999
$ws.Range("M71").Value = -30008.5
$ws.Range("N71").Value = -42487
$ws.Range("H126").Value = 23205.773
$ws.Range("J126").Value = 7149
$ws.Range("L126").Value = 21447
$ws.Range("N126").Value = -26387

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 18352.475
$ws.Range("I107").Value = 3054.2307
$ws.Range("K107").Value = 9162.6921
$ws.Range("M107").Value = -7242.6921
